$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$titleShape = $s.Shapes.Item(3)

# ------------------------------------------------------------------
# Duplicate the title textbox FIRST, while it still has its original
# two runs ("Physical level " / "of abstraction "). That way the new
# shape's two runs - and their distinct rPr (dirty/smtClean) flags -
# come from the real two-run source instead of an already-merged one.
# ------------------------------------------------------------------
$newShape = $titleShape.Duplicate()
$newShape.Name = "TextBox 4"

# ------------------------------------------------------------------
# 1) "Physical level " + "of abstraction " -> single run
#    "Physical level of abstraction "
# ------------------------------------------------------------------
$titleTr = $titleShape.TextFrame.TextRange

# Rewrite the first run (originally "Physical level ", 15 chars) so it
# contains the full merged sentence.
$titleRun1 = $titleTr.Characters(1, 15)
$titleRun1.Text = "Physical level of abstraction "

# The old second run ("of abstraction ", 15 chars) now lives right after
# the freshly-expanded first run; clear it out so only one run remains.
$titleRun2 = $titleTr.Characters(31, 15)
$titleRun2.Text = ""

# Position/size in points; nudged by +0.5 EMU worth of points so the
# float -> EMU conversion lands exactly on the target EMU values
# (1302783, 1917700, 5410200, 646331) instead of rounding one unit low.
$newShape.Left = (1302783 + 0.5) / 12700
$newShape.Top = (1917700 + 0.5) / 12700
$newShape.Width = (5410200 + 0.5) / 12700
$newShape.Height = (646331 + 0.5) / 12700

$bodyTr = $newShape.TextFrame.TextRange

# Reuse the duplicated shape's first run (still "Physical level ", 15
# chars) for the new first sentence...
$bodyRun1 = $bodyTr.Characters(1, 15)
$bodyRun1.Text = "Basically to store data in such a way that it is easy for us to access and fetch "

# ...and its second run (now shifted to start right after run1) for the
# closing sentence.
$bodyRun2 = $bodyTr.Characters(82, 15)
$bodyRun2.Text = "data easily."

Write-Output "ok"
